$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD, AE, AF with the same header style as the
# existing header cells. Copy formatting from AC1 (bold/bordered/
# centered header style) so the new cells reuse the exact same style
# index instead of generating a near-duplicate style.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill team record values for all data rows (2-50).
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 64  # AD = Wins
    $ws.Cells.Item($r, 31).Value = 98  # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF = Ties
}
